$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 4 (Aline Silva / professor de educação física) entirely,
# shifting subsequent rows up.
$ws.Rows.Item(4).Delete()
